$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove existing hyperlinks so they can be re-added in row order with
# sequential relationship ids (rId1..rId21), matching rows 2..22.
$ws.Hyperlinks.Delete()

# Final dataset for rows 2..22: 取得日時, タイトル, カテゴリ, 価格, 締切, URL, 優先度スコア, スキル概要
$data = @(
    @("2025-11-20 12:36:45", "【ECシステム開発】販売データ分析・AI提案・競合監視を統合した販売支援システム構築", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5437916", 410, "🔥AI,Ai ◆開発,システム開発"),
    @("2025-11-20 12:36:45", "【AI開発】生成AI・RAGシステム構築パートナー募集", "システム開発", "1,000 ~ 5,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5437447", 375, "🔥AI,Ai ◆開発"),
    @("2025-11-20 12:36:45", "【急募】Web管理システム構築・AI機能実装のプロを探しています", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5437717", 360, "🔥AI,Ai ◇管理"),
    @("2025-11-20 12:36:45", "大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5427956", 310, "🔥AI,Ai"),
    @("2025-11-20 12:36:45", "【急募】案件管理システム開発のフリーランス募集!", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5437726", 160, "◆開発,システム開発 ◇管理"),
    @("2025-11-20 12:36:45", "React製の予約サイトの調査・不具合修正ができるエンジニアを募集", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5438035", 153, "🔥React ◇サイト"),
    @("2025-11-20 12:36:45", "【急募】Android/iOS対応メンタルヘルスアプリ開発者募集", "システム開発", "1,000,000 円 ~ 3,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5438164", 100, "◆開発 ◇アプリ"),
    @("2025-11-20 12:36:45", "【急募】マーケティングリサーチ用WEBアプリ開発の依頼", "システム開発", "1,000,000 円 ~ 3,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5438081", 100, "◆開発 ◇アプリ"),
    @("2025-11-20 12:36:45", "webアプリの開発", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5437832", 100, "◆開発 ◇アプリ"),
    @("2025-11-20 12:36:45", "信頼される医療コンサル会社のHP作成依頼", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5437868", 55, "◆コンサル"),
    @("2025-11-20 12:36:45", "【急募】価格更新サイトにエクセルアップロード後、内容を更新するプログラム作成依頼", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5437655", 38, "◇サイト"),
    @("2025-11-20 12:36:45", "【急募】WordPressにe-SCOTT決済機能を導入", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5437728", 33, "○WordPress"),
    @("2025-11-20 12:36:45", "エクセルから個々のデータを抽出し、自動でメールを送信するシステム", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5437997", 33, $null),
    @("2025-11-20 12:36:45", "【急募】SNS運用でFX自動売買システムの利用者を増やしたい!", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5437991", 28, $null),
    @("2025-11-20 12:36:45", "【フルスタックエンジニア】 働きながらスキルアップもできるEC業界で活躍してみませんか?", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5429335", 25, $null),
    @("2025-11-20 12:36:45", "【急募】Flutterflowの扱えるノーコードエンジニアを探しています!", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5437783", 25, $null),
    @("2025-11-20 12:36:45", "【SESエンジニア募集】多様なプロジェクトに参画可能!", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5437544", 25, $null),
    @("2025-11-20 12:36:45", "Networkエンジニア", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5432661", 25, $null),
    @("2025-11-20 12:36:45", "AmazonのASINを渡すだけで楽天へ大量出品できる方(数万点規模)】", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5438052", 13, $null),
    @("2025-11-20 12:36:45", "【急募】Wartalesの武器アイコンとモデルを日本刀に差し替え", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5436248", 10, $null),
    @("2025-11-20 12:36:45", "サーバー移管(2ドメイン)のご依頼", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5438014", 10, $null)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]

    $urlCell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($urlCell, $row[5])
    $urlCell.Style = "Hyperlink"
}

# Column width adjustments (B: 46 -> 51, D: 30 -> 32 characters).
# ColumnWidth applies a +5/6 pixel-rounding offset in this engine, so the
# assigned value is pre-compensated to land on the exact target width.
$ws.Columns.Item(2).ColumnWidth = 50.166666666666664
$ws.Columns.Item(4).ColumnWidth = 31.166666666666668
